# mrp_dataset.xlsx edit: add "Analysis" sheet content, tweak Sheet1 data/selection,
# rename Sheet2 -> Analysis, and make Analysis the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: data correction + selection change -----------------------------
# Part 149A7567-911 (row 6) actual demand corrected 181 -> 155; this also
# ripples into the Forecast_Error (H6) / Exception_Flag (I6) formulas.
$ws1.Range("C6").Value = 155

# Selection moves off H6 (no longer tabSelected) down to B36.
$ws1.Range("B36").Select() | Out-Null

# --- Sheet2: build the "Analysis" summary content ---------------------------
$ws2.Range("A2").Value = "Total Parts Analyzed"
$ws2.Range("B2").Formula = "=COUNTA(Sheet1!A2:A21)"

$ws2.Range("A3").Value = "Number of Stockouts"
$ws2.Range("B3").Formula = '=COUNTIF(Sheet1!G2:G21,"Yes")'

$ws2.Range("A4").Value = "Number of Exceptions"
$ws2.Range("B4").Formula = '=COUNTIF(Sheet1!I2:I21,"Exception")'

$ws2.Range("A5").Value = "Average Forecast Error"
$ws2.Range("B5").Formula = "=AVERAGE(Sheet1!H2:H21)"

$ws2.Range("A6").Value = "Max Forecast Error"
$ws2.Range("B6").Formula = "=MAX(Sheet1!H2:H21)"

$ws2.Range("A7").Value = "Min Forecast Error"
$ws2.Range("B7").Formula = "=+MIN(Sheet1!H2:H21)"

$ws2.Range("A10").Value = "Key Observations"

$ws2.Range("A11").Value = "A significant portion of parts triggered exception flags, indicating instability in planning inputs."
$ws2.Range("A11:H11").Merge() | Out-Null
$ws2.Range("A11:H11").HorizontalAlignment = -4108

$ws2.Range("A12").Value = "Stockouts appeared even when forecasts were relatively close, suggesting sensitivity to lead time and safety stock."
$ws2.Range("A12:J12").Merge() | Out-Null
$ws2.Range("A12:J12").HorizontalAlignment = -4108

$ws2.Range("A13").Value = "Exception-based review would allow planners to focus on high-risk parts instead of reviewing everything manually."
$ws2.Range("A13:J13").Merge() | Out-Null
$ws2.Range("A13:J13").HorizontalAlignment = -4108

$ws2.Range("A15").Value = "Why This Matters"

$ws2.Range("A16").Value = "Demonstrates how simple analytics can reveal hidden inefficiencies."
$ws2.Range("A16:F16").Borders.Item(9).LineStyle = 1
$ws2.Range("A16:F16").Borders.Item(9).LineStyle = -4142
$ws2.Range("G16:H16").Borders.Item(9).LineStyle = 1
$ws2.Range("G16:H16").Borders.Item(9).LineStyle = -4142

$ws2.Range("A17").Value = "Shows potential to improve operational efficiency without expensive system changes."

$ws2.Range("A18").Value = "Supports broader application in manufacturing environments facing similar planning challenges."

# Column width for the label column (~26.3 characters to fit the labels).
$ws2.Columns.Item(1).ColumnWidth = 25.5

# A couple of small floating text boxes (artifacts carried over from the
# author's source) anchored near row 10 of the Analysis sheet.
$ws2.Shapes.AddTextbox(1, 31.5, 144.75, 14.55, 20.83) | Out-Null
$ws2.Shapes.AddTextbox(1, 7.5, 141.75, 14.55, 20.83) | Out-Null

# --- Workbook: rename + re-order active tab ---------------------------------
$ws2.Name = "Analysis"
$ws2.Range("A22").Select() | Out-Null
$ws2.Activate()
